# RPA datasets push 2023-11-07
#
# Source data had one extra underwriter row ("상상인" / "상상인제4호스팩")
# that needs to be removed, and the "밀리의서재" row (under "미래") needs
# to move so it sits right after the "두산로보틱스" row within the "미래"
# group instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "상상인" / "상상인제4호스팩" row (row 13). Everything below
#    shifts up by one row.
$ws.Rows(13).Delete()

# 2) Relocate the "밀리의서재" row (now row 8, still above the 두산로보틱스
#    "미래" row that is now row 10) so it comes right after 두산로보틱스.
#    Stage its contents off to the side first (Range.Copy preserves the
#    original cell types/formatting, unlike reading/writing .Value which
#    would coerce the date-looking text back into real dates).
$ws.Range("A8:L8").Copy($ws.Range("A100:L100"))

# Remove the original "밀리의서재" row; rows below shift up again.
$ws.Rows(8).Delete()

# Open a slot right after 두산로보틱스 (which is now row 10) and drop the
# staged row into it.
$ws.Rows(11).Insert()
$ws.Range("A100:L100").Copy($ws.Range("A11:L11"))

# Tidy up the scratch row used for staging.
$ws.Range("A100:L100").Clear()
